# Facility Rental Report - employee role integration and UI updates
# Replaces the three named-facility rows (Function Hall / University Auditorium /
# University Gymnasium) and the old TOTAL row with a single blank "TOTAL" row,
# and clears the three rows below it back to blank template rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 10 becomes the bold "TOTAL" row (was "Function Hall (ACAD Bldg.)" row)
# ---------------------------------------------------------------------------
$a10 = $ws.Range("A10")
$a10.Value = "TOTAL"
$a10.WrapText = $true
$a10.HorizontalAlignment = -4108
$a10.Font.Bold = $true

$b10 = $ws.Range("B10")
$b10.ClearContents()
$b10.Font.Bold = $true

$c10 = $ws.Range("C10")
$c10.Value = 0
$c10.Font.Bold = $true

$d10 = $ws.Range("D10")
$d10.ClearContents()
$d10.Font.Bold = $true

$e10 = $ws.Range("E10")
$e10.Font.Bold = $true
$e10.Value = "'x"
$e10.ClearContents()

# ---------------------------------------------------------------------------
# Row 11 becomes a blank template row (was "University Auditorium" row)
# ---------------------------------------------------------------------------
$ws.Range("A11").ClearContents()
$ws.Range("B11").ClearContents()

$c11 = $ws.Range("C11")
$c11.ClearContents()
$c11.NumberFormat = '_(* #,##0.00_);_(* \(#,##0.00\);_(* "-"??_);_(@_)'

$ws.Range("D11").ClearContents()

$e11 = $ws.Range("E11")
$e11.Value = "'x"
$e11.ClearContents()

# ---------------------------------------------------------------------------
# Row 12 becomes a blank template row (was "University Gymnasium" row)
# ---------------------------------------------------------------------------
$ws.Range("A12").ClearContents()
$ws.Range("B12").ClearContents()

$c12 = $ws.Range("C12")
$c12.ClearContents()
$c12.NumberFormat = '_(* #,##0.00_);_(* \(#,##0.00\);_(* "-"??_);_(@_)'

$ws.Range("D12").ClearContents()

$e12 = $ws.Range("E12")
$e12.Value = "'x"
$e12.ClearContents()

# ---------------------------------------------------------------------------
# Row 13 becomes a blank template row (was the bold "TOTAL" row, now removed)
# ---------------------------------------------------------------------------
$a13 = $ws.Range("A13")
$a13.ClearContents()
$a13.Font.Bold = $false

$b13 = $ws.Range("B13")
$b13.ClearContents()
$b13.Font.Bold = $false

$c13 = $ws.Range("C13")
$c13.ClearContents()
$c13.Font.Bold = $false
$c13.NumberFormat = '_(* #,##0.00_);_(* \(#,##0.00\);_(* "-"??_);_(@_)'

$d13 = $ws.Range("D13")
$d13.ClearContents()
$d13.Font.Bold = $false

$e13 = $ws.Range("E13")
$e13.Font.Bold = $false
$e13.Value = "'x"
$e13.ClearContents()

# ---------------------------------------------------------------------------
# Move the active selection from the old TOTAL row to the new TOTAL row
# ---------------------------------------------------------------------------
$ws.Range("A10:E10").Select() | Out-Null

Write-Host "Facility rental report rows 10-13 updated"
